$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "avgs" label row
$ws.Range("A29").Value = "avgs"

# Averages of the three sample groups (soja S, milho M, mata nativa F)
$ws.Range("B29:C31").NumberFormat = "0.0000"

$ws.Range("B29").Formula = "=AVERAGE(B8:B13)"
$ws.Range("C29").Formula = "=AVERAGE(C8:C13)"

$ws.Range("B30").Formula = "=AVERAGE(B14:B19)"
$ws.Range("C30").Formula = "=AVERAGE(C14:C19)"

$ws.Range("B31").Formula = "=AVERAGE(B20:B25)"
$ws.Range("C31").Formula = "=AVERAGE(C20:C25)"

[void]$ws.Range("A30").Select()
